$wb = $excel.ActiveWorkbook

# "464cc037-93c5-457b-9dc0-def937e079c7.md" has moved from "Ready for handoff" to "In Translation".
# Update the status on every sheet that tracks it (row 5 in each table).

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E5").Value = "In Translation"
$wsOverview.Range("F5").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C5").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C5").Value = "In Translation"
